# Staging.SubOutputPersonRole.xlsx - column header rename
# (SubOutputPersonRole_ID/PersonSourceKey/RoleSourceKey/SubOutputSourceKey
#  -> PersonBusinessKey/RoleBusinessKey/SubOutputBusinessKey/SubOutputPersonRole_ID)
#
# Row 1 (A1) is an untouched "internal use only" note.
# Row 2 holds the four bold/underlined column-header labels in A2:D2 -
# only their text changes; formatting/layout is left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "PersonBusinessKey"
$ws.Range("B2").Value = "RoleBusinessKey"
$ws.Range("C2").Value = "SubOutputBusinessKey"
$ws.Range("D2").Value = "SubOutputPersonRole_ID"
